$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 4, shifting existing rows 4-29 down to 6-31
$ws.Range("A4:A5").EntireRow.Insert()

# Copy formatting (bold, border, alignment) from A3 to the new A4:A5 index cells
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New row 4 (index 2): HKL = "Holden"
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Holden"
$ws.Range("C4").Value = 0.947080866833913
$ws.Range("D4").Value = 0.8830486967638586
$ws.Range("E4").Value = 1.112796908515802
$ws.Range("F4").Value = 0.9346628205061608
$ws.Range("G4").Value = 0.9346628205061608
$ws.Range("H4").Value = 1.220338010373251
$ws.Range("I4").Value = 1.220338010373251
$ws.Range("J4").Value = 0.952191174442774
$ws.Range("K4").Value = 0.9346628205061608
$ws.Range("L4").Value = 0.952191174442774
$ws.Range("M4").Value = 1.086264592408012
$ws.Range("N4").Value = 1.086264592408012
$ws.Range("O4").Value = 1.095108697777275
$ws.Range("P4").Value = 1.035730668440728
$ws.Range("Q4").Value = 1.035730668440728
$ws.Range("R4").Value = 1.010463706457086
$ws.Range("S4").Value = 1.010463706457086
$ws.Range("T4").Value = 1.008353079572627

# New row 5 (index 3): HKL = "Rizzie Spiral"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Rizzie Spiral"
$ws.Range("C5").Value = 1.046940379947563
$ws.Range("D5").Value = 1.226690916419899
$ws.Range("E5").Value = 0.8853566742153822
$ws.Range("F5").Value = 0.9811535867659517
$ws.Range("G5").Value = 0.9811535867659517
$ws.Range("H5").Value = 0.842144732987282
$ws.Range("I5").Value = 0.842144732987282
$ws.Range("J5").Value = 1.07401314151932
$ws.Range("K5").Value = 0.9811535867659517
$ws.Range("L5").Value = 1.07401314151932
$ws.Range("M5").Value = 0.958078937253301
$ws.Range("N5").Value = 0.958078937253301
$ws.Range("O5").Value = 0.9338381829073281
$ws.Range("P5").Value = 0.9657704870908512
$ws.Range("Q5").Value = 0.9657704870908512
$ws.Range("R5").Value = 0.9696162620096264
$ws.Range("S5").Value = 0.9696162620096264
$ws.Range("T5").Value = 1.009383238642567

# Rename "Thomas Hex" -> "Matthies Hex" (now at row 11 after the shift)
$ws.Range("B11").Value = "Matthies Hex"
